$d = $word.ActiveDocument

# Map of (AlternativeText description, current picture "name") -> new "name"
# as recorded on the wp:docPr / pic:cNvPr of each logo picture in the
# document's headers/footers.
#   - BTec_Logo-Orange : image1.jpg -> image2.jpg
#   - PearsonLogo.png  : image2.png -> image1.png   (appears twice)

function Rename-LogoPicture {
    param($range)

    if ($range -eq $null) { return }
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $ishp = $range.InlineShapes($i)
        $desc = $ishp.AlternativeText

        $newName = $null
        if ($desc -eq "BTec_Logo-Orange") {
            $newName = "image2.jpg"
        } elseif ($desc -like "*PearsonLogo.png") {
            $newName = "image1.png"
        }

        if ($newName -ne $null) {
            $shp = $ishp.ConvertToShape()
            $shp.Name = $newName
            $shp.ConvertToInlineShape() | Out-Null
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections($s)

    for ($h = 1; $h -le 3; $h++) {
        Rename-LogoPicture $sec.Headers($h).Range
    }
    for ($f = 1; $f -le 3; $f++) {
        Rename-LogoPicture $sec.Footers($f).Range
    }
}
